# Schedule 관련 Proto 생성 / ProtoGenerater 복합 pk, list 읽기 작업
# Adds a new EScheduleType enum block (NONE / GACHA / ATTENDANCE) to the
# "Common" sheet, right after the existing EFormationPositionType block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Common")

$ws.Range("A46").Value = "EScheduleType"
$ws.Range("B46").Value = "NONE"
$ws.Range("C46").Value = 0

$ws.Range("A47").Value = "EScheduleType"
$ws.Range("B47").Value = "GACHA"
$ws.Range("C47").Value = 10

$ws.Range("A48").Value = "EScheduleType"
$ws.Range("B48").Value = "ATTENDANCE"
$ws.Range("C48").Value = 20

$ws.Range("C47").Select()
